$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap date values between rows 3 and 4 (column D)
$ws.Range("D3").Value = 44981
$ws.Range("D4").Value = 44980

# Swap volume values between rows 3 and 4 (column M)
$ws.Range("M3").Value = 30
$ws.Range("M4").Value = 50
